# Inicio de la interpretacion de la pendiente de una recta como razon
# constante de cambio o velocidad de cambio.
#
# - Adds the "Tarea 4" column header (F1) for the new task.
# - Fixes a typo in a student's name (IBAnEZ -> IBAÑEZ).
# - Sets column widths for the two new columns.
# - Fills in a missing grade (D31).
# - Leaves the active selection on the new header cell (F1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the student's last name.
$ws.Range("A14").Value = "IBAÑEZ MORENO FRANCISCO LUIS"

# New assignment column header.
$ws.Range("F1").Value = "Tarea 4 visualizacion de funciones en excel con video"

# Widen the new columns to fit their content (~14.66 and ~44.07 characters).
$ws.Columns.Item(5).ColumnWidth = 13.85
$ws.Columns.Item(6).ColumnWidth = 43.15

# Fill in the previously-missing grade for USUGA GEORGE JHON ESTIVEN.
$ws.Range("D31").Value = 5

# Leave the selection on the newly added header cell.
$ws.Range("F1").Select()
